$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D2" = "28.119.76"
    "E2" = "  +4.06%  "
    "D3" = "1.804.47"
    "E3" = "  +4.38%  "
    "D4" = "0.9967"
    "E4" = "  -0.51%  "
    "D5" = "316.43"
    "E5" = "  +2.09%  "
    "D6" = "0.9986"
    "E6" = "  -0.27%  "
    "D7" = "0.5707"
    "E7" = "  +17.78%  "
    "D8" = "0.3850"
    "E8" = "  +10.02%  "
    "D9" = "0.07641"
    "E9" = "  +5.26%  "
    "D10" = "43.08"
    "E10" = "  -0.59%  "
    "E11" = "  +8.19%  "
    "B12" = "Solana"
    "C12" = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
    "D12" = "21.25"
    "E12" = "  +6.46%  "
    "B13" = "BinanceUSD"
    "C13" = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
    "D13" = "0.9967"
    "E13" = "  -0.51%  "
    "D14" = "6.227"
    "E14" = "  +5.88%  "
    "D15" = "1.798.96"
    "E15" = "  +4.16%  "
    "D16" = "7.174"
    "E16" = "  +4.42%  "
    "D17" = "91.96"
    "E17" = "  +5.50%  "
    "D18" = "0.00001081"
    "E18" = "  +4.39%  "
    "D19" = "0.06504"
    "E19" = "  +1.56%  "
    "E20" = "  -0.31%  "
    "D21" = "17.23"
    "E21" = "  +3.62%  "
    "D22" = "5.965"
    "E22" = "  +4.53%  "
    "D23" = "28.135.35"
    "E23" = "  +3.87%  "
    "D24" = "11.25"
    "E24" = "  +3.19%  "
    "D25" = "2.095"
    "E25" = "  +0.82%  "
    "D26" = "20.68"
    "E26" = "  +3.74%  "
    "D27" = "156.15"
    "E27" = "  +1.21%  "
    "D28" = "2.370"
    "E28" = "  +14.11%  "
    "D29" = "2.008.06"
    "E29" = "  +4.42%  "
    "D30" = "123.05"
    "E30" = "  +1.74%  "
    "D31" = "1.143"
    "E31" = "  +8.68%  "
    "D32" = "0.1051"
    "E32" = "  +12.37%  "
    "D33" = "5.723"
    "E33" = "  +6.47%  "
    "D34" = "3.626"
    "E34" = "  -0.79%  "
    "D35" = "0.02307"
    "E35" = "  +5.51%  "
    "D36" = "0.2118"
    "E36" = "  +6.11%  "
    "D37" = "8.638"
    "E37" = "  +15.14%  "
    "D38" = "11.59"
    "E38" = "  +5.49%  "
    "D39" = "5.022"
    "E39" = "  +5.35%  "
    "D40" = "0.06050"
    "E40" = "  +1.41%  "
    "D41" = "0.6297"
    "E41" = "  +5.20%  "
    "D42" = "0.9985"
    "E42" = "  -0.24%  "
    "D43" = "1.400"
    "E43" = "  -2.04%  "
    "D44" = "1.151"
    "E44" = "  +4.63%  "
    "D45" = "13.38"
    "E45" = "  +4.21%  "
    "D46" = "0.5915"
    "E46" = "  +5.10%  "
    "D47" = "3.680"
    "E47" = "  +2.83%  "
    "D48" = "121.78"
    "E48" = "  +2.59%  "
    "D49" = "1.931"
    "E49" = "  +4.45%  "
    "D50" = "1.143"
    "E50" = "  +3.72%  "
    "D51" = "0.06800"
    "E51" = "  +2.34%  "
}

$ws.Range("D2:D51").NumberFormat = "@"

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
